# "after hosting week review"
# - Removes the two oldest order rows (previously rows 31 & 32) from the top
#   of the data block. In terms of row *numbers*, this is equivalent to
#   inserting two blank rows right before the first data row (old row 31),
#   which pushes all existing data rows down by 2 (old row 31 -> row 33,
#   ..., old row 57 -> row 59).
# - Appends two brand-new order rows (618667 / 783506) at the bottom
#   (new rows 60 and 61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (old rows 31-57) down by two rows so that
# it starts at row 33, matching the final layout in the diff.
$ws.Rows("31:32").Insert()

# ---------------------------------------------------------------------
# Row 60 - new order "618667"
# ---------------------------------------------------------------------
# Materialize the row first so the later format-paste is not a no-op.
$ws.Range("A60").Value = "PLACEHOLDER"

# Force the Order ID to be stored as text (matching every other Order ID
# cell, which are shared strings even though they look numeric).
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "618667"
$ws.Range("E60").Value = '{"name":"sulu","mobile":9747561989,"homeAddress":"clerus house adimalathura chowara p.o","city":"ATHIYANNUR","postalCode":695501}'
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 22500
$ws.Range("H60").Value = 45233.50763043981
$ws.Range("I60").Value = "Confirmed"
$ws.Range("J60").Value = "Cash on delivery"
$ws.Range("K60").Value = "pending"

# Restore the standard data-row style (font size 12, general format) on
# the whole row now that all the values (and their types) are set; this
# must be the LAST operation on the row so the temporary "@" text format
# style used above is not left applied to A60.
$ws.Range("A59:K59").Copy()
$ws.Range("A60:K60").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------
# Row 61 - new order "783506"
# ---------------------------------------------------------------------
$ws.Range("A61").Value = "PLACEHOLDER"

$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = "783506"
$ws.Range("E61").Value = '{"name":"sulu","mobile":9747561989,"homeAddress":"clerus house adimalathura chowara p.o","city":"ATHIYANNUR","postalCode":695501}'
$ws.Range("F61").Value = 2250
$ws.Range("G61").Value = 6750
$ws.Range("H61").Value = 45234.547509525466
$ws.Range("I61").Value = "Confirmed"
$ws.Range("J61").Value = "Cash on delivery"
$ws.Range("K61").Value = "pending"

$ws.Range("A59:K59").Copy()
$ws.Range("A61:K61").PasteSpecial(-4122) # xlPasteFormats

Write-Host "applied hosting-week-review update"
